# Rodinia GPU benchmark re-run: update column L ("time2" / relative metric)
# for rows 60-79 (Sheet1) with the newly measured values.
# Values are stored as literal text (10 decimal places, incl. trailing
# zeros) -- same representation the sheet already used -- so a leading
# apostrophe forces Excel to keep them as text instead of coercing to a
# number (which would silently drop the trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    60 = "1.1325580000"
    61 = "1.2671750000"
    62 = "0.8565370000"
    63 = "1.9598740000"
    64 = "0.8027820000"
    65 = "1.1964120000"
    66 = "1.2765400000"
    67 = "1.3033110000"
    68 = "0.2607990000"
    69 = "0.3520820000"
    70 = "1.3297930000"
    71 = "1.2440830000"
    72 = "2.9664410000"
    73 = "0.0102400000"
    74 = "0.9495210000"
    75 = "1.6567660000"
    76 = "4.4983170000"
    78 = "4.0395960000"
    79 = "1.2007830000"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 12).Value = "'" + $updates[$row]
}
